$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.370.71"
$ws.Range("E2").Value = "  -1.11%  "

$ws.Range("D3").Value = "2.954.93"
$ws.Range("E3").Value = "  -2.42%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'538.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.01%  "

$ws.Range("D6").Value = "'134.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.20%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").Value = "2.950.23"
$ws.Range("E8").Value = "  -2.50%  "

$ws.Range("D9").Value = "'0.481"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.47%  "

$ws.Range("D10").Value = "'6.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.57%  "

$ws.Range("D11").Value = "'0.145"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.52%  "

$ws.Range("D12").Value = "'0.439"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.56%  "

$ws.Range("D13").Value = "'0.0000215"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.13%  "

$ws.Range("D14").Value = "'33.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.92%  "

$ws.Range("D15").Value = "3.391.96"
$ws.Range("E15").Value = "  -3.42%  "

$ws.Range("D16").Value = "61.274.96"
$ws.Range("E16").Value = "  -1.33%  "

$ws.Range("D17").Value = "'0.107"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.04%  "

$ws.Range("D18").Value = "2.939.53"
$ws.Range("E18").Value = "  -2.65%  "

$ws.Range("D19").Value = "'6.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.56%  "

$ws.Range("D20").Value = "'459.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.77%  "

$ws.Range("D21").Value = "'13.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.54%  "

$ws.Range("D22").Value = "'0.642"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.56%  "

$ws.Range("D23").Value = "'6.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.11%  "

$ws.Range("D24").Value = "'78.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.30%  "

$ws.Range("D25").Value = "'12.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.72%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.26%  "

$ws.Range("D27").Value = "'2.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.39%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'7.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.53%  "

$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("D30").Value = "'1.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.30%  "

$ws.Range("D31").Value = "'24.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.88%  "

$ws.Range("D32").Value = "'1.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.46%  "

$ws.Range("D33").Value = "'2.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.12%  "

$ws.Range("D34").Value = "'5.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.14%  "

$ws.Range("D35").Value = "'53.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.33%  "

$ws.Range("D36").Value = "'5.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.21%  "

$ws.Range("D37").Value = "'441.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.50%  "

$ws.Range("D38").Value = "'0.0789"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.51%  "

$ws.Range("D39").Value = "'0.0382"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.42%  "

$ws.Range("D40").Value = "2.892.71"
$ws.Range("E40").Value = "  -10.14%  "

$ws.Range("D41").Value = "'0.114"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.92%  "

$ws.Range("D42").Value = "'7.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.94%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.53%  "

$ws.Range("D44").Value = "'26.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.69%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").Value = "'0.243"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.12%  "

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'1.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.55%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.107"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("D49").Value = "'113.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.54%  "

$ws.Range("D50").Value = "0.0₃0480"
$ws.Range("E50").Value = "  -3.64%  "

$ws.Range("D51").Value = "'1.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.00%  "
